$d = $word.ActiveDocument

# --- 1. Student name cell: merge "Liucija Paulina " + "Adomaviciute" runs,
#        drop the spell-check proofErr markers around the surname. ---
$find = $d.Content.Find
$find.Execute("Liucija Paulina Adomaviciute", $false, $false, $false, $false, $false, $true, 1, $false, "Liucija Paulina Adomaviciute", 2) | Out-Null

# --- 2. "...calculation results mixed in for every message." -> drop "mixed in" ---
$find = $d.Content.Find
$find.Execute("calculation results mixed in for every message", $false, $false, $false, $false, $false, $true, 1, $false, "calculation results for every message", 2) | Out-Null

# --- 3. "...may be used with another input." -> drop trailing period (next run already has one) ---
$find = $d.Content.Find
$find.Execute("which later may be used with another input.", $false, $false, $false, $false, $false, $true, 1, $false, "which later may be used with another input", 2) | Out-Null

# --- 4. "...the output secrets of the Diffie-Hellman become the inputs to the root chain." ->
#        "...the output of the Diffie-Hellman becomes the input for the root chain." ---
$find = $d.Content.Find
$find.Execute("the output secrets of the Diffie-Hellman become the inputs to the root chain.", $false, $false, $false, $false, $false, $true, 1, $false, "the output of the Diffie-Hellman becomes the input for the root chain.", 2) | Out-Null

# --- 5. "the Ratchet happens in four steps" -> "the ratchet happens in four steps" ---
$find = $d.Content.Find
$find.Execute("Bob, the Ratchet happens in four steps", $false, $false, $false, $false, $false, $true, 1, $false, "Bob, the ratchet happens in four steps", 2) | Out-Null

# --- 6. "This results in receiving a chain key" -> "The result is a receiving a chain key" ---
$find = $d.Content.Find
$find.Execute("This results in receiving a chain key", $false, $false, $false, $false, $false, $true, 1, $false, "The result is a receiving a chain key", 2) | Out-Null

# --- 7/8/9. Drop proofErr spell-check wrappers around BPubK / APrivK / APubK ---
$find = $d.Content.Find
$find.Execute("BPubK – Bob", $false, $false, $false, $false, $false, $true, 1, $false, "BPubK – Bob", 2) | Out-Null
$find = $d.Content.Find
$find.Execute("APrivK – Alice", $false, $false, $false, $false, $false, $true, 1, $false, "APrivK – Alice", 2) | Out-Null
$find = $d.Content.Find
$find.Execute("APubK – Alice", $false, $false, $false, $false, $false, $true, 1, $false, "APubK – Alice", 2) | Out-Null
